$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 211. This shifts the existing rows
# 211..347 down to 212..348, preserving all of their data/formatting.
$ws.Rows.Item(211).Insert()

# Populate the newly inserted row 211 with the new data point.
$ws.Cells.Item(211, 1).Value = 10
$ws.Cells.Item(211, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(211, 3).Value = "La Araucanía"
$ws.Cells.Item(211, 4).Value = 44603
$ws.Cells.Item(211, 5).Value = 9
$ws.Cells.Item(211, 6).Value = 100114014
$ws.Cells.Item(211, 7).Value = "Betarraga"
$ws.Cells.Item(211, 8).Value = "Sin especificar"
$ws.Cells.Item(211, 9).Value = "Primera"
$ws.Cells.Item(211, 10).Value = 65
$ws.Cells.Item(211, 11).Value = 8000
$ws.Cells.Item(211, 12).Value = 8000
$ws.Cells.Item(211, 13).Value = 8000
$ws.Cells.Item(211, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(211, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(211, 16).Value = 320
$ws.Cells.Item(211, 17).Value = 25
$ws.Cells.Item(211, 18).Value = "Hortaliza"
